# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Re-order match-result rows by swapping the full row contents (columns B:AC)
# between rows, while leaving column A (the sequential row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pairwise swap: row 122 <-> row 123 -------------------------------
$rng122 = $ws.Range("B122:AC122")
$rng123 = $ws.Range("B123:AC123")
$v122 = $rng122.Value()
$v123 = $rng123.Value()
$rng122.Value = $v123
$rng123.Value = $v122

# --- Pairwise swap: row 124 <-> row 125 -------------------------------
$rng124 = $ws.Range("B124:AC124")
$rng125 = $ws.Range("B125:AC125")
$v124 = $rng124.Value()
$v125 = $rng125.Value()
$rng124.Value = $v125
$rng125.Value = $v124

# --- 4-way cyclic re-order among rows 134, 135, 136, 137 ---------------
# New row 134 <- old row 137
# New row 135 <- old row 136
# New row 136 <- old row 134
# New row 137 <- old row 135
$rng134 = $ws.Range("B134:AC134")
$rng135 = $ws.Range("B135:AC135")
$rng136 = $ws.Range("B136:AC136")
$rng137 = $ws.Range("B137:AC137")

$v134 = $rng134.Value()
$v135 = $rng135.Value()
$v136 = $rng136.Value()
$v137 = $rng137.Value()

$rng134.Value = $v137
$rng135.Value = $v136
$rng136.Value = $v134
$rng137.Value = $v135
